$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Column C ("Förändrad") holds a date serial that is refreshed by one day
# on every automated run. Update every data row (row 2 through the last
# used row) from 46061 to 46062.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46061) {
        $cell.Value = 46062
    }
}
